# Update Name of Algo
# Applies the updated RandomForest imputation result values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 16.5803
$ws.Range("D7").Value = -6.9125
$ws.Range("C9").Value = -10.20700000000001
$ws.Range("D12").Value = -6.512799999999999
$ws.Range("C13").Value = -12.61519999999999
$ws.Range("D14").Value = -7.269200000000004
$ws.Range("E15").Value = 16.0768
$ws.Range("C16").Value = -12.64469999999999
$ws.Range("C18").Value = -11.5952
$ws.Range("D19").Value = -8.0548
$ws.Range("C20").Value = -11.57660000000002
$ws.Range("C26").Value = -12.5593
$ws.Range("D26").Value = -8.382600000000004
$ws.Range("C27").Value = -12.9466
$ws.Range("D27").Value = -8.860900000000001
$ws.Range("E28").Value = 16.68600000000001
$ws.Range("C29").Value = -13.30509999999999
$ws.Range("D29").Value = -7.990200000000003
$ws.Range("E33").Value = 17.12660000000001
$ws.Range("C35").Value = -11.51040000000001
$ws.Range("E35").Value = 16.62460000000001
$ws.Range("C36").Value = -12.48030000000001
$ws.Range("D37").Value = -7.8412
$ws.Range("D38").Value = -7.2873
$ws.Range("E38").Value = 17.24599999999999
$ws.Range("E43").Value = 17.2991
$ws.Range("E44").Value = 16.8607
$ws.Range("C45").Value = -13.27939999999998
$ws.Range("E45").Value = 16.86510000000001
$ws.Range("D47").Value = -7.630799999999998
$ws.Range("E47").Value = 16.70950000000001
$ws.Range("D51").Value = -8.679399999999998
$ws.Range("E51").Value = 16.4311
$ws.Range("D52").Value = -7.860500000000004
$ws.Range("E54").Value = 16.93150000000001
$ws.Range("C55").Value = -14.11720000000001
$ws.Range("D55").Value = -8.859899999999991
$ws.Range("C57").Value = -13.6709
$ws.Range("E57").Value = 16.2383
$ws.Range("E62").Value = 16.33750000000001
$ws.Range("E63").Value = 18.29170000000002
$ws.Range("E67").Value = 17.28420000000002
$ws.Range("C69").Value = -11.74089999999999
$ws.Range("D69").Value = -7.150199999999995
$ws.Range("D70").Value = -7.252500000000003
$ws.Range("E70").Value = 17.69620000000002
$ws.Range("C76").Value = -12.13480000000001
$ws.Range("D76").Value = -7.706
$ws.Range("C78").Value = -11.37010000000001
$ws.Range("D81").Value = -8.360899999999999
$ws.Range("E81").Value = 16.2974
$ws.Range("C82").Value = -11.4592
$ws.Range("C83").Value = -14.26780000000001
$ws.Range("D83").Value = -8.654899999999998
$ws.Range("E88").Value = 16.3486
$ws.Range("C93").Value = -10.7514
$ws.Range("D94").Value = -7.142299999999999
$ws.Range("E96").Value = 16.15379999999999
$ws.Range("C97").Value = -12.6584
$ws.Range("E99").Value = 17.1588
$ws.Range("D100").Value = -8.275800000000002
$ws.Range("D102").Value = -7.827499999999998
